$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fruta / hortaliza, semanal -- refresh Femacal de La Calera / Mandarina rows
# 410-422 with the latest weekly prices and append the new rows 423-428.
$rowData = @{
    410 = @{ "A"=3; "B"="Femacal de La Calera"; "C"="Coquimbo"; "D"=44448; "E"=5; "F"="Fruta"; "G"=100102; "H"="Cítricos"; "I"=100102004; "J"="Mandarina"; "K"="Murcott"; "L"="Especial"; "M"=90; "N"=6000; "O"=6000; "P"=6000; "Q"="$/bandeja 10 kilos"; "R"="Provincia de Quillota"; "S"=600; "T"=10 }
    411 = @{ "A"=3; "B"="Femacal de La Calera"; "C"="Coquimbo"; "D"=44448; "E"=5; "F"="Fruta"; "G"=100102; "H"="Cítricos"; "I"=100102004; "J"="Mandarina"; "K"="Murcott"; "L"="Especial"; "M"=70; "N"=8000; "O"=8000; "P"=8000; "Q"="$/caja 15 kilos"; "R"="Provincia de Quillota"; "S"=533; "T"=15 }
    412 = @{ "A"=3; "B"="Femacal de La Calera"; "C"="Coquimbo"; "D"=44448; "E"=5; "F"="Fruta"; "G"=100102; "H"="Cítricos"; "I"=100102004; "J"="Mandarina"; "K"="Murcott"; "L"="Primera"; "M"=89; "N"=5000; "O"=5000; "P"=5000; "Q"="$/bandeja 10 kilos"; "R"="Provincia de Quillota"; "S"=500; "T"=10 }
    413 = @{ "A"=3; "B"="Femacal de La Calera"; "C"="Coquimbo"; "D"=44448; "E"=5; "F"="Fruta"; "G"=100102; "H"="Cítricos"; "I"=100102004; "J"="Mandarina"; "K"="Murcott"; "L"="Primera"; "M"=80; "N"=7000; "O"=7000; "P"=7000; "Q"="$/caja 15 kilos"; "R"="Provincia de Quillota"; "S"=467; "T"=15 }
    414 = @{ "A"=3; "B"="Femacal de La Calera"; "C"="Coquimbo"; "D"=44448; "E"=5; "F"="Fruta"; "G"=100102; "H"="Cítricos"; "I"=100102004; "J"="Mandarina"; "K"="Murcott"; "L"="Segunda"; "M"=90; "N"=4000; "O"=4000; "P"=4000; "Q"="$/bandeja 10 kilos"; "R"="Provincia de Quillota"; "S"=400; "T"=10 }
    415 = @{ "A"=3; "B"="Femacal de La Calera"; "C"="Coquimbo"; "D"=44448; "E"=5; "F"="Fruta"; "G"=100102; "H"="Cítricos"; "I"=100102004; "J"="Mandarina"; "K"="Murcott"; "L"="Segunda"; "M"=67; "N"=6000; "O"=6000; "P"=6000; "Q"="$/caja 15 kilos"; "R"="Provincia de Quillota"; "S"=400; "T"=15 }
    416 = @{ "A"=3; "B"="Femacal de La Calera"; "C"="Coquimbo"; "D"=44238; "E"=5; "F"="Fruta"; "G"=100102; "H"="Cítricos"; "I"=100102004; "J"="Mandarina"; "K"="Clemenuless"; "L"="Primera"; "M"=90; "N"=7000; "O"=7000; "P"=7000; "Q"="$/bandeja 10 kilos"; "R"="Provincia de Quillota"; "S"=700; "T"=10 }
    417 = @{ "A"=3; "B"="Femacal de La Calera"; "C"="Coquimbo"; "D"=44399; "E"=5; "F"="Fruta"; "G"=100102; "H"="Cítricos"; "I"=100102004; "J"="Mandarina"; "K"="Clemenuless"; "L"="Especial"; "M"=85; "N"=7000; "O"=7000; "P"=7000; "Q"="$/bandeja 10 kilos"; "R"="Provincia de Quillota"; "S"=700; "T"=10 }
    418 = @{ "A"=3; "B"="Femacal de La Calera"; "C"="Coquimbo"; "D"=44399; "E"=5; "F"="Fruta"; "G"=100102; "H"="Cítricos"; "I"=100102004; "J"="Mandarina"; "K"="Clemenuless"; "L"="Primera"; "M"=87; "N"=6000; "O"=6000; "P"=6000; "Q"="$/bandeja 10 kilos"; "R"="Provincia de Quillota"; "S"=600; "T"=10 }
    419 = @{ "A"=3; "B"="Femacal de La Calera"; "C"="Coquimbo"; "D"=44399; "E"=5; "F"="Fruta"; "G"=100102; "H"="Cítricos"; "I"=100102004; "J"="Mandarina"; "K"="Clemenuless"; "L"="Segunda"; "M"=80; "N"=5000; "O"=5000; "P"=5000; "Q"="$/bandeja 10 kilos"; "R"="Provincia de Quillota"; "S"=500; "T"=10 }
    420 = @{ "A"=3; "B"="Femacal de La Calera"; "C"="Coquimbo"; "D"=44399; "E"=5; "F"="Fruta"; "G"=100102; "H"="Cítricos"; "I"=100102004; "J"="Mandarina"; "K"="Murcott"; "L"="Especial"; "M"=65; "N"=8000; "O"=8000; "P"=8000; "Q"="$/bandeja 10 kilos"; "R"="Provincia de San Felipe de Aconcagua"; "S"=800; "T"=10 }
    421 = @{ "A"=3; "B"="Femacal de La Calera"; "C"="Coquimbo"; "D"=44399; "E"=5; "F"="Fruta"; "G"=100102; "H"="Cítricos"; "I"=100102004; "J"="Mandarina"; "K"="Murcott"; "L"="Primera"; "M"=70; "N"=7000; "O"=7000; "P"=7000; "Q"="$/bandeja 10 kilos"; "R"="Provincia de San Felipe de Aconcagua"; "S"=700; "T"=10 }
    422 = @{ "A"=3; "B"="Femacal de La Calera"; "C"="Coquimbo"; "D"=44399; "E"=5; "F"="Fruta"; "G"=100102; "H"="Cítricos"; "I"=100102004; "J"="Mandarina"; "K"="Murcott"; "L"="Segunda"; "M"=70; "N"=6000; "O"=6000; "P"=6000; "Q"="$/bandeja 10 kilos"; "R"="Provincia de San Felipe de Aconcagua"; "S"=600; "T"=10 }
    423 = @{ "A"=3; "B"="Femacal de La Calera"; "C"="Coquimbo"; "D"=44400; "E"=5; "F"="Fruta"; "G"=100102; "H"="Cítricos"; "I"=100102004; "J"="Mandarina"; "K"="Clemenuless"; "L"="Especial"; "M"=85; "N"=7000; "O"=7000; "P"=7000; "Q"="$/bandeja 10 kilos"; "R"="Provincia de Quillota"; "S"=700; "T"=10 }
    424 = @{ "A"=3; "B"="Femacal de La Calera"; "C"="Coquimbo"; "D"=44400; "E"=5; "F"="Fruta"; "G"=100102; "H"="Cítricos"; "I"=100102004; "J"="Mandarina"; "K"="Clemenuless"; "L"="Primera"; "M"=87; "N"=6000; "O"=6000; "P"=6000; "Q"="$/bandeja 10 kilos"; "R"="Provincia de Quillota"; "S"=600; "T"=10 }
    425 = @{ "A"=3; "B"="Femacal de La Calera"; "C"="Coquimbo"; "D"=44400; "E"=5; "F"="Fruta"; "G"=100102; "H"="Cítricos"; "I"=100102004; "J"="Mandarina"; "K"="Clemenuless"; "L"="Segunda"; "M"=80; "N"=5000; "O"=5000; "P"=5000; "Q"="$/bandeja 10 kilos"; "R"="Provincia de Quillota"; "S"=500; "T"=10 }
    426 = @{ "A"=3; "B"="Femacal de La Calera"; "C"="Coquimbo"; "D"=44400; "E"=5; "F"="Fruta"; "G"=100102; "H"="Cítricos"; "I"=100102004; "J"="Mandarina"; "K"="Murcott"; "L"="Especial"; "M"=56; "N"=8000; "O"=8000; "P"=8000; "Q"="$/bandeja 10 kilos"; "R"="Provincia de San Felipe de Aconcagua"; "S"=800; "T"=10 }
    427 = @{ "A"=3; "B"="Femacal de La Calera"; "C"="Coquimbo"; "D"=44400; "E"=5; "F"="Fruta"; "G"=100102; "H"="Cítricos"; "I"=100102004; "J"="Mandarina"; "K"="Murcott"; "L"="Primera"; "M"=67; "N"=7000; "O"=7000; "P"=7000; "Q"="$/bandeja 10 kilos"; "R"="Provincia de San Felipe de Aconcagua"; "S"=700; "T"=10 }
    428 = @{ "A"=3; "B"="Femacal de La Calera"; "C"="Coquimbo"; "D"=44400; "E"=5; "F"="Fruta"; "G"=100102; "H"="Cítricos"; "I"=100102004; "J"="Mandarina"; "K"="Murcott"; "L"="Segunda"; "M"=60; "N"=6000; "O"=6000; "P"=6000; "Q"="$/bandeja 10 kilos"; "R"="Provincia de San Felipe de Aconcagua"; "S"=600; "T"=10 }
}

foreach ($r in $rowData.Keys) {
    $row = $rowData[$r]
    foreach ($col in $row.Keys) {
        $ws.Range("$col$r").Value = $row[$col]
    }
}

# Rows 423-428 are brand new; give their Fecha (D) cells the same
# date number format already used by the rest of the column.
foreach ($r in 423..428) {
    $ws.Range("D$r").NumberFormat = $ws.Range("D410").NumberFormat
}
